$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    "Tử Vi đồng cung Thiên Cơ tại cung đối Nô Bộc",
    "Tử Vi đồng cung Thái Dương tại cung đối Nô Bộc",
    "Tử Vi đồng cung Vũ Khúc tại cung đối Nô Bộc",
    "Tử Vi đồng cung Thiên Đồng tại cung đối Nô Bộc",
    "Tử Vi đồng cung Liêm Trinh tại cung đối Nô Bộc",
    "Tử Vi đồng cung Thiên Phủ tại cung đối Nô Bộc",
    "Tử Vi đồng cung Thái Âm tại cung đối Nô Bộc",
    "Tử Vi đồng cung Tham Lang tại cung đối Nô Bộc",
    "Tử Vi đồng cung Cự Môn tại cung đối Nô Bộc",
    "Tử Vi đồng cung Thiên Tướng tại cung đối Nô Bộc",
    "Tử Vi đồng cung Thiên Lương tại cung đối Nô Bộc",
    "Tử Vi đồng cung Thất Sát tại cung đối Nô Bộc",
    "Tử Vi đồng cung Phá Quân tại cung đối Nô Bộc",
    "Thiên Cơ đồng cung Thái Dương tại cung đối Nô Bộc",
    "Thiên Cơ đồng cung Vũ Khúc tại cung đối Nô Bộc",
    "Thiên Cơ đồng cung Thiên Đồng tại cung đối Nô Bộc",
    "Thiên Cơ đồng cung Liêm Trinh tại cung đối Nô Bộc",
    "Thiên Cơ đồng cung Thiên Phủ tại cung đối Nô Bộc",
    "Thiên Cơ đồng cung Thái Âm tại cung đối Nô Bộc",
    "Thiên Cơ đồng cung Tham Lang tại cung đối Nô Bộc",
    "Thiên Cơ đồng cung Cự Môn tại cung đối Nô Bộc",
    "Thiên Cơ đồng cung Thiên Tướng tại cung đối Nô Bộc",
    "Thiên Cơ đồng cung Thiên Lương tại cung đối Nô Bộc",
    "Thiên Cơ đồng cung Thất Sát tại cung đối Nô Bộc",
    "Thiên Cơ đồng cung Phá Quân tại cung đối Nô Bộc",
    "Thái Dương đồng cung Vũ Khúc tại cung đối Nô Bộc",
    "Thái Dương đồng cung Thiên Đồng tại cung đối Nô Bộc",
    "Thái Dương đồng cung Liêm Trinh tại cung đối Nô Bộc",
    "Thái Dương đồng cung Thiên Phủ tại cung đối Nô Bộc",
    "Thái Dương đồng cung Thái Âm tại cung đối Nô Bộc",
    "Thái Dương đồng cung Tham Lang tại cung đối Nô Bộc",
    "Thái Dương đồng cung Cự Môn tại cung đối Nô Bộc",
    "Thái Dương đồng cung Thiên Tướng tại cung đối Nô Bộc",
    "Thái Dương đồng cung Thiên Lương tại cung đối Nô Bộc",
    "Thái Dương đồng cung Thất Sát tại cung đối Nô Bộc",
    "Thái Dương đồng cung Phá Quân tại cung đối Nô Bộc",
    "Vũ Khúc đồng cung Thiên Đồng tại cung đối Nô Bộc",
    "Vũ Khúc đồng cung Liêm Trinh tại cung đối Nô Bộc",
    "Vũ Khúc đồng cung Thiên Phủ tại cung đối Nô Bộc",
    "Vũ Khúc đồng cung Thái Âm tại cung đối Nô Bộc",
    "Vũ Khúc đồng cung Tham Lang tại cung đối Nô Bộc",
    "Vũ Khúc đồng cung Cự Môn tại cung đối Nô Bộc",
    "Vũ Khúc đồng cung Thiên Tướng tại cung đối Nô Bộc",
    "Vũ Khúc đồng cung Thiên Lương tại cung đối Nô Bộc",
    "Vũ Khúc đồng cung Thất Sát tại cung đối Nô Bộc",
    "Vũ Khúc đồng cung Phá Quân tại cung đối Nô Bộc",
    "Thiên Đồng đồng cung Liêm Trinh tại cung đối Nô Bộc",
    "Thiên Đồng đồng cung Thiên Phủ tại cung đối Nô Bộc",
    "Thiên Đồng đồng cung Thái Âm tại cung đối Nô Bộc",
    "Thiên Đồng đồng cung Tham Lang tại cung đối Nô Bộc",
    "Thiên Đồng đồng cung Cự Môn tại cung đối Nô Bộc",
    "Thiên Đồng đồng cung Thiên Tướng tại cung đối Nô Bộc",
    "Thiên Đồng đồng cung Thiên Lương tại cung đối Nô Bộc",
    "Thiên Đồng đồng cung Thất Sát tại cung đối Nô Bộc",
    "Thiên Đồng đồng cung Phá Quân tại cung đối Nô Bộc",
    "Liêm Trinh đồng cung Thiên Phủ tại cung đối Nô Bộc",
    "Liêm Trinh đồng cung Thái Âm tại cung đối Nô Bộc",
    "Liêm Trinh đồng cung Tham Lang tại cung đối Nô Bộc",
    "Liêm Trinh đồng cung Cự Môn tại cung đối Nô Bộc",
    "Liêm Trinh đồng cung Thiên Tướng tại cung đối Nô Bộc",
    "Liêm Trinh đồng cung Thiên Lương tại cung đối Nô Bộc",
    "Liêm Trinh đồng cung Thất Sát tại cung đối Nô Bộc",
    "Liêm Trinh đồng cung Phá Quân tại cung đối Nô Bộc",
    "Thiên Phủ đồng cung Thái Âm tại cung đối Nô Bộc",
    "Thiên Phủ đồng cung Tham Lang tại cung đối Nô Bộc",
    "Thiên Phủ đồng cung Cự Môn tại cung đối Nô Bộc",
    "Thiên Phủ đồng cung Thiên Tướng tại cung đối Nô Bộc",
    "Thiên Phủ đồng cung Thiên Lương tại cung đối Nô Bộc",
    "Thiên Phủ đồng cung Thất Sát tại cung đối Nô Bộc",
    "Thiên Phủ đồng cung Phá Quân tại cung đối Nô Bộc",
    "Thái Âm đồng cung Tham Lang tại cung đối Nô Bộc",
    "Thái Âm đồng cung Cự Môn tại cung đối Nô Bộc",
    "Thái Âm đồng cung Thiên Tướng tại cung đối Nô Bộc",
    "Thái Âm đồng cung Thiên Lương tại cung đối Nô Bộc",
    "Thái Âm đồng cung Thất Sát tại cung đối Nô Bộc",
    "Thái Âm đồng cung Phá Quân tại cung đối Nô Bộc",
    "Tham Lang đồng cung Cự Môn tại cung đối Nô Bộc",
    "Tham Lang đồng cung Thiên Tướng tại cung đối Nô Bộc",
    "Tham Lang đồng cung Thiên Lương tại cung đối Nô Bộc",
    "Tham Lang đồng cung Thất Sát tại cung đối Nô Bộc",
    "Tham Lang đồng cung Phá Quân tại cung đối Nô Bộc",
    "Cự Môn đồng cung Thiên Tướng tại cung đối Nô Bộc",
    "Cự Môn đồng cung Thiên Lương tại cung đối Nô Bộc",
    "Cự Môn đồng cung Thất Sát tại cung đối Nô Bộc",
    "Cự Môn đồng cung Phá Quân tại cung đối Nô Bộc",
    "Thiên Tướng đồng cung Thiên Lương tại cung đối Nô Bộc",
    "Thiên Tướng đồng cung Thất Sát tại cung đối Nô Bộc",
    "Thiên Tướng đồng cung Phá Quân tại cung đối Nô Bộc",
    "Thiên Lương đồng cung Thất Sát tại cung đối Nô Bộc",
    "Thiên Lương đồng cung Phá Quân tại cung đối Nô Bộc",
    "Thất Sát đồng cung Phá Quân tại cung đối Nô Bộc",
    "Tử Vi tọa thủ tại cung đối Nô Bộc",
    "Thiên Cơ tọa thủ tại cung đối Nô Bộc",
    "Thái Dương tọa thủ tại cung đối Nô Bộc",
    "Vũ Khúc tọa thủ tại cung đối Nô Bộc",
    "Thiên Đồng tọa thủ tại cung đối Nô Bộc",
    "Liêm Trinh tọa thủ tại cung đối Nô Bộc",
    "Thiên Phủ tọa thủ tại cung đối Nô Bộc",
    "Thái Âm tọa thủ tại cung đối Nô Bộc",
    "Tham Lang tọa thủ tại cung đối Nô Bộc",
    "Cự Môn tọa thủ tại cung đối Nô Bộc",
    "Thiên Tướng tọa thủ tại cung đối Nô Bộc",
    "Thiên Lương tọa thủ tại cung đối Nô Bộc",
    "Thất Sát tọa thủ tại cung đối Nô Bộc",
    "Phá Quân tọa thủ tại cung đối Nô Bộc"
)

$startRow = 4300
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $values[$i]
    $ws.Cells.Item($row, 2).Value = $values[$i]
}

# Match the author's final selection/view state (last block of new rows).
$ws.Range("B4391:B4404").Select()

